$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new transaction row (row 26) to the trade log, mirroring the
# data/format of the existing rows (the table previously ended at row 25).
$row = 26

$ws.Cells.Item($row, 1).Value = 35                      # A: S.No.
$ws.Cells.Item($row, 2).Value = 42257                   # B: Date
$ws.Cells.Item($row, 3).Value = "Fund Sundaram Select"  # C: Name
$ws.Cells.Item($row, 4).Value = 3422                    # D: Transaction Detail
$ws.Cells.Item($row, 5).Value = "paid"                  # E: Transaction Tax
$ws.Cells.Item($row, 6).Value = 95.88                   # F: Quantity
$ws.Cells.Item($row, 8).Value = 237900                  # H: Sale

# Carry over the number formats used by row 25 for the date and quantity
# columns (the other columns use the default/general format already).
$ws.Range("B25").Copy()
$ws.Range("B26").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F25").Copy()
$ws.Range("F26").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H27").Select()
